# Netzreglerberechnungen für Leitungen hinzugefügt
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: drop the "[in %]" suffix, now expressed as a ratio instead of a percentage ---
$ws.Range("E1").Value = "aktuelle Leistung pL "

# --- Update the line values: D (Bemessungsleistung) and E (aktuelle Leistung, now a fraction) ---
$ws.Range("D2").Value = 95
$ws.Range("E2").Value = 1

$ws.Range("E3").Value = 0.5

$ws.Range("E4").Value = -0.7

$ws.Range("E5").Value = 0.2

$ws.Range("E6").Value = 0.1

$ws.Range("E7").Value = -0.6

$ws.Range("D8").Value = 300
$ws.Range("E8").Value = 0.3

# --- New "Netzregler" calculation rows, centered, below the table ---
$ws.Range("D16:F16").HorizontalAlignment = -4108

$ws.Range("D28:F28").HorizontalAlignment = -4108
$ws.Range("D28:F28").Font.Bold = $true

# --- Restore the active selection to the last-edited cell ---
$ws.Range("E8").Select() | Out-Null
